# Generate Report for Handoff
#
# Marks the rows that are "Ready for handoff" (and were last generated at
# the previous run's timestamp) as handed off: their Priority is stamped
# "ht" and their handoff-generation timestamps are bumped to the new run
# time, on the Overview sheet as well as each per-locale sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Rows (on every sheet) that belong to files just handed off in this run.
$rows = @(7, 8, 9, 11, 13, 14)

foreach ($r in $rows) {
    # Overview sheet: column G = "Latest HO Xliff Generate Date"
    $overview.Range("G$r").Value = "2016-09-06 20:28:19"

    # zh-cn sheet: column E = "Priority", column H = "Latest Handoff Datetime"
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-09-06 20:28:04"

    # de-de sheet: column E = "Priority", column H = "Latest Handoff Datetime"
    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-09-06 20:28:19"
}
